$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.026972333333333
$ws.Range("H2").Value = 3.080917
$ws.Range("I2").Value = 0.2032541865322035
$ws.Range("J2").Value = 0.2032541865322035
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 118.0470123333333
$ws.Range("N2").Value = 354.141037
$ws.Range("O2").Value = 0.4657216250363638
$ws.Range("P2").Value = 0.4657216250363638
$ws.Range("Q2").Value = 121.2310156989921
$ws.Range("R2").Value = 1091.079141290929
$ws.Range("S2").Value = 0.09465987004722201
$ws.Range("T2").Value = 0.09465987004722203

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.026972333333333
$ws.Range("H3").Value = 3.080917
$ws.Range("I3").Value = 0.2032541865322035
$ws.Range("J3").Value = 0.2032541865322035
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 59.36586533333332
$ws.Range("N3").Value = 178.097596
$ws.Range("O3").Value = 0.2342114953037475
$ws.Range("P3").Value = 0.2342114953037476
$ws.Range("Q3").Value = 60.96710124172577
$ws.Range("R3").Value = 548.703911175532
$ws.Range("S3").Value = 0.0476044669544542
$ws.Range("T3").Value = 0.04760446695445422

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.026972333333333
$ws.Range("H4").Value = 3.080917
$ws.Range("I4").Value = 0.2032541865322035
$ws.Range("J4").Value = 0.2032541865322035
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 56.84506433333333
$ws.Range("N4").Value = 170.535193
$ws.Range("O4").Value = 0.2242663767030476
$ws.Range("P4").Value = 0.2242663767030477
$ws.Range("Q4").Value = 58.37830835688678
$ws.Range("R4").Value = 525.4047752119809
$ws.Range("S4").Value = 0.04558307996330266
$ws.Range("T4").Value = 0.04558307996330267

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.026972333333333
$ws.Range("H5").Value = 3.080917
$ws.Range("I5").Value = 0.2032541865322035
$ws.Range("J5").Value = 0.2032541865322035
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.21324333333333
$ws.Range("N5").Value = 57.63973
$ws.Range("O5").Value = 0.07580050295684103
$ws.Range("P5").Value = 0.07580050295684104
$ws.Range("Q5").Value = 19.73146933693445
$ws.Range("R5").Value = 177.58322403241
$ws.Range("S5").Value = 0.01540676956722461
$ws.Range("T5").Value = 0.01540676956722461

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.358031333333334
$ws.Range("H6").Value = 4.074094000000001
$ws.Range("I6").Value = 0.2687760370778347
$ws.Range("J6").Value = 0.2687760370778347
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 118.0470123333333
$ws.Range("N6").Value = 354.141037
$ws.Range("O6").Value = 0.4657216250363638
$ws.Range("P6").Value = 0.4657216250363638
$ws.Range("Q6").Value = 160.3115415550531
$ws.Range("R6").Value = 1442.803873995478
$ws.Range("S6").Value = 0.1251748127587231
$ws.Range("T6").Value = 0.1251748127587231

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.358031333333334
$ws.Range("H7").Value = 4.074094000000001
$ws.Range("I7").Value = 0.2687760370778347
$ws.Range("J7").Value = 0.2687760370778347
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 59.36586533333332
$ws.Range("N7").Value = 178.097596
$ws.Range("O7").Value = 0.2342114953037475
$ws.Range("P7").Value = 0.2342114953037476
$ws.Range("Q7").Value = 80.62070525311378
$ws.Range("R7").Value = 725.586347278024
$ws.Range("S7").Value = 0.06295043754581514
$ws.Range("T7").Value = 0.06295043754581515

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.358031333333334
$ws.Range("H8").Value = 4.074094000000001
$ws.Range("I8").Value = 0.2687760370778347
$ws.Range("J8").Value = 0.2687760370778347
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 56.84506433333333
$ws.Range("N8").Value = 170.535193
$ws.Range("O8").Value = 0.2242663767030476
$ws.Range("P8").Value = 0.2242663767030477
$ws.Range("Q8").Value = 77.1973785100158
$ws.Range("R8").Value = 694.7764065901421
$ws.Range("S8").Value = 0.06027742798004997
$ws.Range("T8").Value = 0.06027742798004997

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.358031333333334
$ws.Range("H9").Value = 4.074094000000001
$ws.Range("I9").Value = 0.2687760370778347
$ws.Range("J9").Value = 0.2687760370778347
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.21324333333333
$ws.Range("N9").Value = 57.63973
$ws.Range("O9").Value = 0.07580050295684103
$ws.Range("P9").Value = 0.07580050295684104
$ws.Range("Q9").Value = 26.09218646162445
$ws.Range("R9").Value = 234.82967815462
$ws.Range("S9").Value = 0.02037335879324642
$ws.Range("T9").Value = 0.02037335879324642

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.667646666666666
$ws.Range("H10").Value = 8.002939999999999
$ws.Range("I10").Value = 0.5279697763899619
$ws.Range("J10").Value = 0.5279697763899619
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 118.0470123333333
$ws.Range("N10").Value = 354.141037
$ws.Range("O10").Value = 0.4657216250363638
$ws.Range("P10").Value = 0.4657216250363638
$ws.Range("Q10").Value = 314.9077189609755
$ws.Range("R10").Value = 2834.169470648779
$ws.Range("S10").Value = 0.2458869422304187
$ws.Range("T10").Value = 0.2458869422304187

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.667646666666666
$ws.Range("H11").Value = 8.002939999999999
$ws.Range("I11").Value = 0.5279697763899619
$ws.Range("J11").Value = 0.5279697763899619
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 59.36586533333332
$ws.Range("N11").Value = 178.097596
$ws.Range("O11").Value = 0.2342114953037475
$ws.Range("P11").Value = 0.2342114953037476
$ws.Range("Q11").Value = 158.3671527702489
$ws.Range("R11").Value = 1425.30437493224
$ws.Range("S11").Value = 0.1236565908034782
$ws.Range("T11").Value = 0.1236565908034782

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.667646666666666
$ws.Range("H12").Value = 8.002939999999999
$ws.Range("I12").Value = 0.5279697763899619
$ws.Range("J12").Value = 0.5279697763899619
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 56.84506433333333
$ws.Range("N12").Value = 170.535193
$ws.Range("O12").Value = 0.2242663767030476
$ws.Range("P12").Value = 0.2242663767030477
$ws.Range("Q12").Value = 151.6425463852689
$ws.Range("R12").Value = 1364.78291746742
$ws.Range("S12").Value = 0.118405868759695
$ws.Range("T12").Value = 0.118405868759695

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.667646666666666
$ws.Range("H13").Value = 8.002939999999999
$ws.Range("I13").Value = 0.5279697763899619
$ws.Range("J13").Value = 0.5279697763899619
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.21324333333333
$ws.Range("N13").Value = 57.63973
$ws.Range("O13").Value = 0.07580050295684103
$ws.Range("P13").Value = 0.07580050295684104
$ws.Range("Q13").Value = 51.25414453402222
$ws.Range("R13").Value = 461.2873008061999
$ws.Range("S13").Value = 0.04002037459637001
$ws.Range("T13").Value = 0.04002037459637001
